$d = $word.ActiveDocument

# Position a collapsed range at the very end of the document body (after the
# last paragraph mark, immediately before the section properties).
$endPos = $d.Content.End
$rng = $d.Range($endPos, $endPos)

# New blank paragraph (matches the standalone <w:p/> in the diff).
$rng.InsertParagraphAfter()
$rng.Start = $rng.End + 1
$rng.End = $rng.Start

# New paragraph holding the "Delete the cluster" instructions.
$rng.InsertParagraphAfter()
$rng.Start = $rng.End + 1
$rng.End = $rng.Start

$rng.InsertAfter("Delete the  cluster:-")
$rng.Collapse(0)

$rng.InsertBreak(6)
$rng.Start = $rng.Start + 1
$rng.End = $rng.Start

$rng.InsertAfter("eksctl delete cluster \")
$rng.Collapse(0)

$rng.InsertBreak(6)
$rng.Start = $rng.Start + 1
$rng.End = $rng.Start

$rng.InsertAfter("--name idli \")
$rng.Collapse(0)

$rng.InsertBreak(6)
$rng.Start = $rng.Start + 1
$rng.End = $rng.Start

$rng.InsertAfter("--region us-east-2")
